$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 38, pushing the current rows 38-57 down to 40-59.
$ws.Range("A38:A39").EntireRow.Insert()

# Row 38: new Cereza / Lapins / Primera entry (week of 2021-12-13, serial 44529)
$ws.Cells.Item(38,1).Value = 7
$ws.Cells.Item(38,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38,3).Value = "Ñuble"
$ws.Cells.Item(38,4).Value = 44529
$ws.Cells.Item(38,5).Value = 16
$ws.Cells.Item(38,6).Value = "Fruta"
$ws.Cells.Item(38,7).Value = 100103
$ws.Cells.Item(38,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(38,9).Value = 100103001
$ws.Cells.Item(38,10).Value = "Cereza"
$ws.Cells.Item(38,11).Value = "Lapins"
$ws.Cells.Item(38,12).Value = "Primera"
$ws.Cells.Item(38,13).Value = 120
$ws.Cells.Item(38,14).Value = 16000
$ws.Cells.Item(38,15).Value = 17000
$ws.Cells.Item(38,16).Value = 16500
$ws.Cells.Item(38,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(38,18).Value = "Provincia de Curicó"
$ws.Cells.Item(38,19).Value = 1650
$ws.Cells.Item(38,20).Value = 10

# Row 39: new Cereza / Lapins / Segunda entry (same date)
$ws.Cells.Item(39,1).Value = 7
$ws.Cells.Item(39,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39,3).Value = "Ñuble"
$ws.Cells.Item(39,4).Value = 44529
$ws.Cells.Item(39,5).Value = 16
$ws.Cells.Item(39,6).Value = "Fruta"
$ws.Cells.Item(39,7).Value = 100103
$ws.Cells.Item(39,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39,9).Value = 100103001
$ws.Cells.Item(39,10).Value = "Cereza"
$ws.Cells.Item(39,11).Value = "Lapins"
$ws.Cells.Item(39,12).Value = "Segunda"
$ws.Cells.Item(39,13).Value = 120
$ws.Cells.Item(39,14).Value = 14000
$ws.Cells.Item(39,15).Value = 15000
$ws.Cells.Item(39,16).Value = 14500
$ws.Cells.Item(39,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(39,18).Value = "Provincia de Curicó"
$ws.Cells.Item(39,19).Value = 1450
$ws.Cells.Item(39,20).Value = 10

# Make sure the date cells use the same date-time number format as the rest of column D.
$ws.Cells.Item(38,4).NumberFormat = $ws.Cells.Item(40,4).NumberFormat
$ws.Cells.Item(39,4).NumberFormat = $ws.Cells.Item(40,4).NumberFormat
